$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (far outside the used A:E range) used to force plain
# numeric-looking strings (e.g. "186.85") to be written as literal text,
# matching the source data where every Price/Volume cell is stored as
# text. Formatted as Text ("@") so Value assignment is not re-parsed as
# a number, then copy/pasted as values into the target cell.
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "74.005.75"
$ws.Range("E2").Value = "  +7.61%  "

# Row 3
$ws.Range("D3").Value = "2.624.97"
$ws.Range("E3").Value = "  +7.48%  "

# Row 4
$helper.Value = "0.999"
$helper.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$helper.Value = "186.85"
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +14.23%  "

# Row 6
$helper.Value = "581.81"
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +3.77%  "

# Row 7
$helper.Value = "0.999"
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -0.15%  "

# Row 8
$helper.Value = "0.531"
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  +4.62%  "

# Row 9
$helper.Value = "0.199"
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +17.33%  "

# Row 10
$ws.Range("D10").Value = "2.620.85"
$ws.Range("E10").Value = "  +7.38%  "

# Row 12
$helper.Value = "0.358"
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +7.77%  "

# Row 13
$helper.Value = "4.67"
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +1.54%  "

# Row 14
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$helper.Value = "0.0000189"
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +5.76%  "

# Row 15
$ws.Range("D15").Value = "73.909.31"
$ws.Range("E15").Value = "  +7.62%  "

# Row 16
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.103.06"
$ws.Range("E16").Value = "  +7.39%  "

# Row 17
$helper.Value = "26.38"
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  +12.66%  "

# Row 18
$ws.Range("D18").Value = "2.628.52"
$ws.Range("E18").Value = "  +7.63%  "

# Row 19
$helper.Value = "9.07"
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +29.22%  "

# Row 20
$helper.Value = "11.83"
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +11.68%  "

# Row 21
$helper.Value = "366.98"
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +8.21%  "

# Row 22
$helper.Value = "2.30"
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +17.90%  "

# Row 23
$helper.Value = "4.07"
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +5.86%  "

# Row 24
$ws.Range("E24").Value = "  -0.11%  "

# Row 25
$helper.Value = "69.79"
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +6.49%  "

# Row 26
$helper.Value = "4.13"
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  +8.86%  "

# Row 27
$helper.Value = "9.34"
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +11.47%  "

# Row 28
$ws.Range("D28").Value = "2.756.01"
$ws.Range("E28").Value = "  +7.35%  "

# Row 29
$helper.Value = "0.999"
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -1.29%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0942"
$ws.Range("E30").Value = "  +14.37%  "

# Row 31
$helper.Value = "521.93"
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +20.34%  "

# Row 32
$helper.Value = "1.38"
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +15.38%  "

# Row 33
$helper.Value = "7.66"
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  +6.55%  "

# Row 34
$helper.Value = "1.75"
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +8.93%  "

# Row 35
$ws.Range("E35").Value = "  -0.04%  "

# Row 36
$helper.Value = "162.51"
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +1.96%  "

# Row 37
$helper.Value = "0.118"
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +10.25%  "

# Row 38
$helper.Value = "19.13"
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +6.26%  "

# Row 39
$helper.Value = "19.27"
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +1.42%  "

# Row 40
$ws.Range("E40").Value = "  +0.05%  "

# Row 41
$helper.Value = "4.90"
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +11.95%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$helper.Value = "1.66"
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  +9.59%  "

# Row 43
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$helper.Value = "0.325"
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  +8.18%  "

# Row 44
$helper.Value = "161.87"
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +24.42%  "

# Row 45
$helper.Value = "2.38"
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  +14.39%  "

# Row 46
$helper.Value = "1.18"
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +9.31%  "

# Row 47
$helper.Value = "38.90"
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  +3.73%  "

# Row 48
$helper.Value = "0.0853"
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +18.49%  "

# Row 49
$helper.Value = "3.61"
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +8.29%  "

# Row 50
$helper.Value = "0.524"
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +8.45%  "

# Row 51
$helper.Value = "20.72"
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +22.07%  "

# Remove the helper column entirely (not just clear it) so no empty
# cell / style trace is left behind and the sheet's used range stays A1:E51.
$helper.EntireColumn.Delete()
